$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Consolidate the "Absent" column (H) so it correctly reflects attendance:
# Absent = 1 - Real (column E), for the data rows in this report.
$ws.Range("H3").Value = 1
$ws.Range("H6").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("H13").Value = 0
